$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.31719433333333
$ws.Range("H2").Value = 120.951583
$ws.Range("I2").Value = 0.1054453461914981
$ws.Range("J2").Value = 0.1054453461914981
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 6189.859227633701
$ws.Range("R2").Value = 55708.73304870331
$ws.Range("S2").Value = 0.03344944620181483
$ws.Range("T2").Value = 0.03344944620181484

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.31719433333333
$ws.Range("H3").Value = 120.951583
$ws.Range("I3").Value = 0.1054453461914981
$ws.Range("J3").Value = 0.1054453461914981
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 6805.530415820886
$ws.Range("R3").Value = 61249.77374238797
$ws.Range("S3").Value = 0.03677647829251835
$ws.Range("T3").Value = 0.03677647829251836

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.31719433333333
$ws.Range("H4").Value = 120.951583
$ws.Range("I4").Value = 0.1054453461914981
$ws.Range("J4").Value = 0.1054453461914981
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 2745.210798049501
$ws.Range("R4").Value = 24706.89718244551
$ws.Range("S4").Value = 0.0148348738678992
$ws.Range("T4").Value = 0.01483487386789921

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 40.31719433333333
$ws.Range("H5").Value = 120.951583
$ws.Range("I5").Value = 0.1054453461914981
$ws.Range("J5").Value = 0.1054453461914981
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 3772.18446968712
$ws.Range("R5").Value = 33949.66022718407
$ws.Range("S5").Value = 0.02038454782926568
$ws.Range("T5").Value = 0.02038454782926568

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 137.4747006666667
$ws.Range("H6").Value = 412.424102
$ws.Range("I6").Value = 0.3595505005759843
$ws.Range("J6").Value = 0.3595505005759843
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 21106.35570154748
$ws.Range("R6").Value = 189957.2013139273
$ws.Range("S6").Value = 0.1140568603569314
$ws.Range("T6").Value = 0.1140568603569314

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 137.4747006666667
$ws.Range("H7").Value = 412.424102
$ws.Range("I7").Value = 0.3595505005759843
$ws.Range("J7").Value = 0.3595505005759843
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 23205.68859672234
$ws.Range("R7").Value = 208851.1973705011
$ws.Range("S7").Value = 0.1254014677469279
$ws.Range("T7").Value = 0.1254014677469279

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 137.4747006666667
$ws.Range("H8").Value = 412.424102
$ws.Range("I8").Value = 0.3595505005759843
$ws.Range("J8").Value = 0.3595505005759843
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 9360.696818546548
$ws.Range("R8").Value = 84246.27136691894
$ws.Range("S8").Value = 0.05058436922856641
$ws.Range("T8").Value = 0.05058436922856642

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 137.4747006666667
$ws.Range("H9").Value = 412.424102
$ws.Range("I9").Value = 0.3595505005759843
$ws.Range("J9").Value = 0.3595505005759843
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 12862.50046424822
$ws.Range("R9").Value = 115762.5041782339
$ws.Range("S9").Value = 0.06950780324355861
$ws.Range("T9").Value = 0.06950780324355861

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 81.80342533333332
$ws.Range("H10").Value = 245.410276
$ws.Range("I10").Value = 0.2139481838098067
$ws.Range("J10").Value = 0.2139481838098067
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 12559.19950592737
$ws.Range("R10").Value = 113032.7955533463
$ws.Range("S10").Value = 0.06786879196475276
$ws.Range("T10").Value = 0.06786879196475279

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 81.80342533333332
$ws.Range("H11").Value = 245.410276
$ws.Range("I11").Value = 0.2139481838098067
$ws.Range("J11").Value = 0.2139481838098067
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 13808.39387338153
$ws.Range("R11").Value = 124275.5448604338
$ws.Range("S11").Value = 0.07461932671087848
$ws.Range("T11").Value = 0.07461932671087849

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 81.80342533333332
$ws.Range("H12").Value = 245.410276
$ws.Range("I12").Value = 0.2139481838098067
$ws.Range("J12").Value = 0.2139481838098067
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 5570.021680720857
$ws.Range("R12").Value = 50130.19512648772
$ws.Range("S12").Value = 0.03009989948082226
$ws.Range("T12").Value = 0.03009989948082227

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 81.80342533333332
$ws.Range("H13").Value = 245.410276
$ws.Range("I13").Value = 0.2139481838098067
$ws.Range("J13").Value = 0.2139481838098067
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 7653.747134742582
$ws.Range("R13").Value = 68883.72421268324
$ws.Range("S13").Value = 0.04136016565335314
$ws.Range("T13").Value = 0.04136016565335315

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 122.756256
$ws.Range("H14").Value = 368.268768
$ws.Range("I14").Value = 0.321055969422711
$ws.Range("J14").Value = 0.321055969422711
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 18846.6473551991
$ws.Range("R14").Value = 169619.8261967919
$ws.Range("S14").Value = 0.1018455983583499
$ws.Range("T14").Value = 0.10184559835835

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 122.756256
$ws.Range("H15").Value = 368.268768
$ws.Range("I15").Value = 0.321055969422711
$ws.Range("J15").Value = 0.321055969422711
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 20721.21951327322
$ws.Range("R15").Value = 186490.975619459
$ws.Range("S15").Value = 0.1119756188074403
$ws.Range("T15").Value = 0.1119756188074403

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 122.756256
$ws.Range("H16").Value = 368.268768
$ws.Range("I16").Value = 0.321055969422711
$ws.Range("J16").Value = 0.321055969422711
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 8358.513162229439
$ws.Range("R16").Value = 75226.61846006497
$ws.Range("S16").Value = 0.04516865829500251
$ws.Range("T16").Value = 0.04516865829500252

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 122.756256
$ws.Range("H17").Value = 368.268768
$ws.Range("I17").Value = 0.321055969422711
$ws.Range("J17").Value = 0.321055969422711
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 11485.40343883229
$ws.Range("R17").Value = 103368.6309494906
$ws.Range("S17").Value = 0.0620660939619182
$ws.Range("T17").Value = 0.0620660939619182

